$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.356.40"
$ws.Range("E2").Value = "  -3.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.499.49"
$ws.Range("E3").Value = "  -4.75%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.11"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.31"
$ws.Range("E6").Value = "  -6.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.498.36"
$ws.Range("E7").Value = "  -4.69%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -4.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("E11").Value = "  -3.77%  "

$ws.Range("E12").Value = "  -4.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -4.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.091.22"
$ws.Range("E14").Value = "  -4.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.44"
$ws.Range("E15").Value = "  -3.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.498.80"
$ws.Range("E16").Value = "  -5.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.259.69"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  -2.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.97"
$ws.Range("E20").Value = "  -5.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.88"
$ws.Range("E21").Value = "  -5.07%  "

$ws.Range("E22").Value = "  -12.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.618"
$ws.Range("E23").Value = "  -4.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.20"
$ws.Range("E24").Value = "  -2.83%  "

$ws.Range("E25").Value = "  +4.83%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.639.93"
$ws.Range("E27").Value = "  -4.72%  "

$ws.Range("E29").Value = "  -5.51%  "

$ws.Range("E30").Value = "  -4.62%  "

$ws.Range("E31").Value = "  -7.40%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.65"
$ws.Range("E34").Value = "  -3.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.12"
$ws.Range("E35").Value = "  -4.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  -6.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.488.45"
$ws.Range("E37").Value = "  -5.04%  "

$ws.Range("E38").Value = "  -3.80%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.47"
$ws.Range("E42").Value = "  -2.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0873"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("E44").Value = "  -6.09%  "

$ws.Range("E45").Value = "  -4.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.43"
$ws.Range("E46").Value = "  -2.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.60"
$ws.Range("E47").Value = "  -5.66%  "

$ws.Range("E48").Value = "  +5.56%  "

$ws.Range("E49").Value = "  -5.81%  "

$ws.Range("E50").Value = "  -4.23%  "

$ws.Range("E51").Value = "  -3.72%  "
